$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was regenerated with new TPM values. Rows 2-4 are updated
# in place with the refreshed numbers (and a couple of cluster-label cells
# change), and the now-obsolete rows 5-7 are removed entirely.

# Row 2 (Sending cluster ECs -> Target cluster FAPs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 34.36677233333334
$ws.Range("H2").Value = 103.100317
$ws.Range("I2").Value = 0.241469236623423
$ws.Range("J2").Value = 0.241469236623423
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4055383333333333
$ws.Range("N2").Value = 1.216615
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 13.93704357410611
$ws.Range("R2").Value = 125.433392166955
$ws.Range("S2").Value = 0.241469236623423
$ws.Range("T2").Value = 0.241469236623423

# Row 3 (Sending cluster FAPs -> FAPs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 72.00815333333334
$ws.Range("H3").Value = 216.02446
$ws.Range("I3").Value = 0.5059466640455351
$ws.Range("J3").Value = 0.5059466640455351
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 29.20206648921111
$ws.Range("R3").Value = 262.8185984029
$ws.Range("S3").Value = 0.5059466640455351
$ws.Range("T3").Value = 0.5059466640455351

# Row 4 now carries the (formerly row 6) MuSCs -> FAPs data
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 35.94867966666666
$ws.Range("H4").Value = 107.846039
$ws.Range("I4").Value = 0.2525840993310418
$ws.Range("J4").Value = 0.2525840993310418
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4055383333333333
$ws.Range("N4").Value = 1.216615
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 14.57856763755389
$ws.Range("R4").Value = 131.207108737985
$ws.Range("S4").Value = 0.2525840993310418
$ws.Range("T4").Value = 0.2525840993310418

# The old rows 5-7 (duplicate MuSCs/FAPs target-cluster permutations) no
# longer exist in the refreshed dataset - remove them entirely.
$ws.Range("A5:A7").EntireRow.Delete()
